# Refresh "previsao_retorno" snapshot (dados bibi e add):
# 1) bump every cached "meses sem comprar" (situacao) bucket by 0.1 month
#    for the rows whose inactivity text was recomputed, and
# 2) rewrite the handful of client rows whose probability / history /
#    date-window metrics shifted in this data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$situacaoUpdates = @{
    3 = "INATIVO - 46.3 meses sem comprar"
    5 = "INATIVO - 46.1 meses sem comprar"
    7 = "INATIVO - 45.9 meses sem comprar"
    8 = "INATIVO - 39.4 meses sem comprar"
    9 = "INATIVO - 39.8 meses sem comprar"
    10 = "INATIVO - 42.9 meses sem comprar"
    11 = "INATIVO - 50.6 meses sem comprar"
    12 = "INATIVO - 46.4 meses sem comprar"
    13 = "INATIVO - 30.3 meses sem comprar"
    14 = "INATIVO - 40.4 meses sem comprar"
    15 = "INATIVO - 49.1 meses sem comprar"
    16 = "INATIVO - 46.2 meses sem comprar"
    17 = "INATIVO - 45.9 meses sem comprar"
    19 = "INATIVO - 49.3 meses sem comprar"
    20 = "INATIVO - 21.4 meses sem comprar"
    21 = "INATIVO - 49.7 meses sem comprar"
    22 = "INATIVO - 12.2 meses sem comprar"
    23 = "INATIVO - 0.6 meses sem comprar"
    24 = "INATIVO - 19.3 meses sem comprar"
    25 = "INATIVO - 34.2 meses sem comprar"
    27 = "INATIVO - 45.5 meses sem comprar"
    28 = "INATIVO - 44.3 meses sem comprar"
    30 = "INATIVO - 47.3 meses sem comprar"
    35 = "INATIVO - 12.6 meses sem comprar"
    36 = "INATIVO - 8.6 meses sem comprar"
    37 = "INATIVO - 28.9 meses sem comprar"
    38 = "INATIVO - 40.3 meses sem comprar"
    42 = "INATIVO - 45.5 meses sem comprar"
    43 = "INATIVO - 36.3 meses sem comprar"
    45 = "INATIVO - 41.7 meses sem comprar"
    47 = "INATIVO - 42.2 meses sem comprar"
    52 = "INATIVO - 22.8 meses sem comprar"
    55 = "INATIVO - 39.0 meses sem comprar"
    56 = "INATIVO - 15.1 meses sem comprar"
    58 = "INATIVO - 0.3 meses sem comprar"
    59 = "INATIVO - 39.6 meses sem comprar"
    60 = "INATIVO - 40.7 meses sem comprar"
    61 = "INATIVO - 39.2 meses sem comprar"
    62 = "INATIVO - 39.2 meses sem comprar"
    63 = "INATIVO - 40.7 meses sem comprar"
    64 = "INATIVO - 39.2 meses sem comprar"
    65 = "INATIVO - 39.2 meses sem comprar"
    66 = "INATIVO - 39.2 meses sem comprar"
    67 = "INATIVO - 43.6 meses sem comprar"
    69 = "INATIVO - 39.2 meses sem comprar"
    71 = "INATIVO - 8.4 meses sem comprar"
    73 = "INATIVO - 40.7 meses sem comprar"
    74 = "INATIVO - 39.6 meses sem comprar"
    75 = "INATIVO - 8.4 meses sem comprar"
    79 = "INATIVO - 39.2 meses sem comprar"
    80 = "INATIVO - 16.1 meses sem comprar"
    82 = "INATIVO - 28.4 meses sem comprar"
    85 = "INATIVO - 17.0 meses sem comprar"
    87 = "INATIVO - 8.6 meses sem comprar"
    89 = "INATIVO - 48.9 meses sem comprar"
    90 = "INATIVO - 25.5 meses sem comprar"
    92 = "INATIVO - 38.8 meses sem comprar"
    93 = "INATIVO - 43.1 meses sem comprar"
    94 = "INATIVO - 39.6 meses sem comprar"
    96 = "INATIVO - 18.8 meses sem comprar"
    99 = "INATIVO - 27.1 meses sem comprar"
    101 = "INATIVO - 40.3 meses sem comprar"
    102 = "INATIVO - 21.6 meses sem comprar"
    103 = "INATIVO - 6.2 meses sem comprar"
    104 = "INATIVO - 38.7 meses sem comprar"
    107 = "INATIVO - 15.8 meses sem comprar"
    108 = "INATIVO - 45.2 meses sem comprar"
    110 = "INATIVO - 41.4 meses sem comprar"
    111 = "INATIVO - 35.9 meses sem comprar"
    112 = "INATIVO - 39.9 meses sem comprar"
    114 = "INATIVO - 30.3 meses sem comprar"
    115 = "INATIVO - 1.9 meses sem comprar"
    116 = "INATIVO - 21.2 meses sem comprar"
    117 = "INATIVO - 25.8 meses sem comprar"
    118 = "INATIVO - 45.5 meses sem comprar"
    119 = "INATIVO - 28.4 meses sem comprar"
    120 = "INATIVO - 27.3 meses sem comprar"
    122 = "INATIVO - 10.9 meses sem comprar"
    125 = "INATIVO - 39.2 meses sem comprar"
    126 = "INATIVO - 28.9 meses sem comprar"
    128 = "INATIVO - 6.4 meses sem comprar"
    129 = "INATIVO - 28.5 meses sem comprar"
    130 = "INATIVO - 36.6 meses sem comprar"
    132 = "INATIVO - 8.3 meses sem comprar"
    136 = "INATIVO - 9.0 meses sem comprar"
    137 = "INATIVO - 34.3 meses sem comprar"
    139 = "INATIVO - 41.3 meses sem comprar"
    140 = "INATIVO - 24.2 meses sem comprar"
    141 = "INATIVO - 20.0 meses sem comprar"
    142 = "INATIVO - 34.9 meses sem comprar"
    144 = "INATIVO - 12.7 meses sem comprar"
    147 = "INATIVO - 10.3 meses sem comprar"
    153 = "INATIVO - 28.4 meses sem comprar"
    156 = "INATIVO - 26.6 meses sem comprar"
    164 = "INATIVO - 14.0 meses sem comprar"
    165 = "INATIVO - 7.2 meses sem comprar"
    167 = "INATIVO - 21.3 meses sem comprar"
    168 = "INATIVO - 13.2 meses sem comprar"
    175 = "INATIVO - 34.5 meses sem comprar"
    184 = "INATIVO - 20.5 meses sem comprar"
    190 = "INATIVO - 19.1 meses sem comprar"
    191 = "INATIVO - 33.8 meses sem comprar"
    192 = "INATIVO - 24.6 meses sem comprar"
    198 = "INATIVO - 39.0 meses sem comprar"
    202 = "INATIVO - 30.5 meses sem comprar"
    208 = "INATIVO - 27.7 meses sem comprar"
    217 = "INATIVO - 30.9 meses sem comprar"
    219 = "INATIVO - 17.3 meses sem comprar"
    229 = "INATIVO - 35.2 meses sem comprar"
    232 = "INATIVO - 23.2 meses sem comprar"
    236 = "INATIVO - 24.9 meses sem comprar"
    243 = "INATIVO - 10.0 meses sem comprar"
    246 = "INATIVO - 14.2 meses sem comprar"
    248 = "INATIVO - 15.2 meses sem comprar"
    252 = "INATIVO - 18.6 meses sem comprar"
    253 = "INATIVO - 28.1 meses sem comprar"
    257 = "INATIVO - 8.0 meses sem comprar"
    260 = "INATIVO - 17.9 meses sem comprar"
    264 = "INATIVO - 24.2 meses sem comprar"
    270 = "INATIVO - 13.9 meses sem comprar"
    274 = "INATIVO - 14.6 meses sem comprar"
    283 = "INATIVO - 15.5 meses sem comprar"
    289 = "INATIVO - 12.0 meses sem comprar"
    294 = "INATIVO - 8.5 meses sem comprar"
    298 = "INATIVO - 10.4 meses sem comprar"
    307 = "INATIVO - 12.0 meses sem comprar"
    312 = "INATIVO - 14.0 meses sem comprar"
    316 = "INATIVO - 8.8 meses sem comprar"
    324 = "INATIVO - 12.2 meses sem comprar"
    327 = "INATIVO - 13.4 meses sem comprar"
    330 = "INATIVO - 7.9 meses sem comprar"
    331 = "INATIVO - 8.1 meses sem comprar"
    333 = "INATIVO - 11.6 meses sem comprar"
    335 = "INATIVO - 9.3 meses sem comprar"
    336 = "INATIVO - 12.7 meses sem comprar"
    340 = "INATIVO - 12.2 meses sem comprar"
    347 = "INATIVO - 11.2 meses sem comprar"
    353 = "INATIVO - 8.1 meses sem comprar"
    405 = "INATIVO - 9.1 meses sem comprar"
}
foreach ($row in $situacaoUpdates.Keys) {
    $ws.Cells.Item($row, 10).Value = $situacaoUpdates[$row]
}

# Per-row metric corrections
# Row 143
$ws.Cells.Item(143, 2).Value = 0.5
$ws.Cells.Item(143, 4).Value = 0.83
$ws.Cells.Item(143, 5).Value = 39
$ws.Cells.Item(143, 6).Value = 0.83
$ws.Cells.Item(143, 8).Value = 45835
$ws.Cells.Item(143, 9).Value = 45896

# Row 189
$ws.Cells.Item(189, 5).Value = 77
$ws.Cells.Item(189, 8).Value = 45835
$ws.Cells.Item(189, 9).Value = 45850

# Row 222
$ws.Cells.Item(222, 2).Value = 0.5
$ws.Cells.Item(222, 4).Value = 0.67
$ws.Cells.Item(222, 5).Value = 12
$ws.Cells.Item(222, 6).Value = 0.67
$ws.Cells.Item(222, 7).Value = "1x a cada 3 meses - irregular (preferencialmente na 2ª quinzena)"
$ws.Cells.Item(222, 8).Value = 45831
$ws.Cells.Item(222, 9).Value = 45923

# Row 223
$ws.Cells.Item(223, 2).Value = 0.58
$ws.Cells.Item(223, 4).Value = 0.83
$ws.Cells.Item(223, 5).Value = 15
$ws.Cells.Item(223, 6).Value = 0.83
$ws.Cells.Item(223, 7).Value = "1x a cada 3 meses - irregular (preferencialmente na 2ª quinzena)"
$ws.Cells.Item(223, 8).Value = 45831
$ws.Cells.Item(223, 9).Value = 45923

# Row 225
$ws.Cells.Item(225, 3).Value = 0.5
$ws.Cells.Item(225, 4).Value = 0.83
$ws.Cells.Item(225, 5).Value = 33
$ws.Cells.Item(225, 6).Value = 0.83
$ws.Cells.Item(225, 8).Value = 45835
$ws.Cells.Item(225, 9).Value = 45896

# Row 239
$ws.Cells.Item(239, 2).Value = 0.58
$ws.Cells.Item(239, 3).Value = 0.33
$ws.Cells.Item(239, 5).Value = 33
$ws.Cells.Item(239, 8).Value = 45835
$ws.Cells.Item(239, 9).Value = 45865

# Row 261
$ws.Cells.Item(261, 5).Value = 44
$ws.Cells.Item(261, 8).Value = 45835
$ws.Cells.Item(261, 9).Value = 45865

# Row 358
$ws.Cells.Item(358, 2).Value = 0.5
$ws.Cells.Item(358, 3).Value = 0.5
$ws.Cells.Item(358, 5).Value = 17
$ws.Cells.Item(358, 8).Value = 45835
$ws.Cells.Item(358, 9).Value = 45865

# Row 368
$ws.Cells.Item(368, 2).Value = 0.58
$ws.Cells.Item(368, 3).Value = 0.5
$ws.Cells.Item(368, 5).Value = 20
$ws.Cells.Item(368, 8).Value = 45835
$ws.Cells.Item(368, 9).Value = 45865

# Row 369
$ws.Cells.Item(369, 2).Value = 0.67
$ws.Cells.Item(369, 3).Value = 0.33
$ws.Cells.Item(369, 5).Value = 19
$ws.Cells.Item(369, 8).Value = 45834
$ws.Cells.Item(369, 9).Value = 45864

# Row 390
$ws.Cells.Item(390, 2).Value = 0.33
$ws.Cells.Item(390, 4).Value = 0.5
$ws.Cells.Item(390, 5).Value = 9
$ws.Cells.Item(390, 6).Value = 0.5
$ws.Cells.Item(390, 8).Value = 45835
$ws.Cells.Item(390, 9).Value = 45896

# Row 395
$ws.Cells.Item(395, 5).Value = 24
$ws.Cells.Item(395, 8).Value = 45835
$ws.Cells.Item(395, 9).Value = 45850

